$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-16 (columns A, B, C)
$data = @(
    @("<paragraph>", "<poil>", 40),
    @("<zero>", "<zero>", 37),
    @("<many>", "<her>", 36),
    @("<sentence>", "<sentence>", 39),
    @("<on>", "<are>", 41),
    @("<the>", "<the>", 37),
    @("<now>", "<now>", 34),
    @("<be>", "<be>", 32),
    @("<some>", "<time>", 42),
    @("<get>", "<an>", 37),
    @("<way>", "<may>", 36),
    @("<tango>", "<in>", 41),
    @("<make>", "<make>", 36),
    @("<up>", "<of>", 40),
    @("<some>", "<down>", 21)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

# Remove old rows 17 and 18 which no longer exist in the target range (A1:C16)
$ws.Range("A17:C18").Delete()
